# Underline the leading "Rate of change" phrase (without its trailing
# space) in the "Rate of change <...>" paragraphs on slides 13, 14 and 15.
# Setting Font.Underline on a sub-range automatically splits the original
# run into "Rate of change" (u="sng") + " " (unchanged) runs, matching the
# target OOXML.

$p = $ppt.ActivePresentation

# Map: slide index -> shape index of the "Content Placeholder 2" shape
# that holds the "Rate of change ..." text.
$targets = @{
    13 = 2
    14 = 4
    15 = 2
}

foreach ($slideIdx in $targets.Keys) {
    $s = $p.Slides.Item($slideIdx)
    $sh = $s.Shapes.Item($targets[$slideIdx])
    $tr = $sh.TextFrame.TextRange

    # "Rate of change" is the first 14 characters of the paragraph.
    $chars = $tr.Characters(1, 14)
    $chars.Font.Underline = 1
}
